$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.017.89"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").Value = "3.387.38"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'571.49"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").Value = "'141.63"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'7.65"
$ws.Range("E9").Value = "  +2.19%  "

$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").Value = "3.968.09"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("E13").Value = "  +1.91%  "

$ws.Range("D14").Value = "'27.81"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").Value = "3.392.96"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").Value = "61.120.91"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "'6.10"
$ws.Range("E18").Value = "  -2.66%  "

$ws.Range("D19").Value = "'13.62"
$ws.Range("E19").Value = "  -3.47%  "

$ws.Range("D20").Value = "'8.89"
$ws.Range("E20").Value = "  -2.31%  "

$ws.Range("D21").Value = "'383.57"
$ws.Range("E21").Value = "  -1.51%  "

$ws.Range("D22").Value = "'75.38"
$ws.Range("E22").Value = "  +3.00%  "

$ws.Range("D23").Value = "'0.551"
$ws.Range("E23").Value = "  -1.68%  "

$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "'0.0000115"
$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("D26").Value = "3.524.68"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  +3.30%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "'7.22"
$ws.Range("E29").Value = "  -2.50%  "

$ws.Range("D30").Value = "'7.97"
$ws.Range("E30").Value = "  -1.31%  "

$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -3.63%  "

$ws.Range("D34").Value = "'23.20"
$ws.Range("E34").Value = "  -2.37%  "

$ws.Range("D35").Value = "'6.94"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").Value = "'165.96"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").Value = "3.422.60"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").Value = "'4.97"
$ws.Range("E38").Value = "  -0.39%  "

$ws.Range("E39").Value = "  -2.76%  "

$ws.Range("D40").Value = "'0.0765"
$ws.Range("E40").Value = "  -1.64%  "

$ws.Range("D41").Value = "'26.80"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.780"
$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").Value = "'1.65"
$ws.Range("E45").Value = "  -2.31%  "

$ws.Range("D46").Value = "'1.12"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").Value = "2.449.53"
$ws.Range("E47").Value = "  -3.44%  "

$ws.Range("D48").Value = "'22.91"
$ws.Range("E48").Value = "  -0.77%  "

$ws.Range("D49").Value = "'6.67"
$ws.Range("E49").Value = "  -2.75%  "

$ws.Range("D50").Value = "'2.14"
$ws.Range("E50").Value = "  +9.77%  "

$ws.Range("E51").Value = "  -2.94%  "
